# Split the single "Tabelle1" timesheet worksheet into two sheets:
#   - "Logging"   : a small config/state sheet (carryover row/column pointers)
#   - "Timesheet" : the original timesheet content (copied from Tabelle1)
#
# The original sheet keeps its identity (sheetId) and becomes "Logging" with
# brand-new content; a duplicate of its original content becomes the new
# "Timesheet" sheet placed right after it.

$wb = $excel.ActiveWorkbook
$orig = $wb.Worksheets.Item(1)

# 1) Duplicate the original timesheet sheet; the copy keeps all the
#    timesheet rows/formulas/formatting and becomes "Timesheet".
$orig.Copy([System.Type]::Missing, $orig)
$timesheet = $wb.Worksheets.Item(2)
$timesheet.Name = "Timesheet"

# 2) Turn the original sheet into the new "Logging" sheet: wipe its old
#    timesheet content/formatting (including merged cells). Note: the
#    logo picture lives on a drawing part shared with "Timesheet", so it
#    is intentionally left in place rather than deleted here.
$orig.Cells.UnMerge()
$orig.Cells.Clear()
$orig.Name = "Logging"

# Small config table used to persist where the timesheet carryover value
# comes from: B1 = label, A2/B2 = row pointer + value, A3/B3 = column
# pointer + value.
$orig.Range("B1").Value = "carryover"
$orig.Range("A2").Value = "row"
$orig.Range("B2").Value = 34
$orig.Range("A3").Value = "column"
$orig.Range("B3").Value = 10
$orig.Range("B2").Select()

# Logging is a plain data sheet - no special page fit-to-page behaviour.
$orig.PageSetup.FitToPage = $false

# Narrow the old D12:G32 block selection down to just D12, and leave the
# Timesheet tab as the active one (matches the authored workbook state).
$timesheet.Activate()
$timesheet.Range("D12").Select()

Write-Host "Sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Host (" - " + $s.Name)
}
